$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "アルマジロ" post row (row 246). Deleting the entire row shifts
# every row below it (247-299) up by one, matching the new A1:C298 extent.
$ws.Rows.Item(246).Delete()
